$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.712771666666666
$ws.Range("H2").Value = 17.138315
$ws.Range("I2").Value = 0.1683613830606884
$ws.Range("J2").Value = 0.1683613830606885
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03970866666666666
$ws.Range("N2").Value = 0.119126
$ws.Range("Q2").Value = 0.2268465458544444
$ws.Range("R2").Value = 2.04161891269
$ws.Range("S2").Value = 0.1683613830606884
$ws.Range("T2").Value = 0.1683613830606885

# Row 3
$ws.Range("I3").Value = 0.4370667227533506
$ws.Range("J3").Value = 0.4370667227533506
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03970866666666666
$ws.Range("N3").Value = 0.119126
$ws.Range("Q3").Value = 0.5888944041804444
$ws.Range("R3").Value = 5.300049637623999
$ws.Range("S3").Value = 0.4370667227533506
$ws.Range("T3").Value = 0.4370667227533506

# Row 4
$ws.Range("G4").Value = 13.388457
$ws.Range("H4").Value = 40.165371
$ws.Range("I4").Value = 0.3945718941859609
$ws.Range("J4").Value = 0.3945718941859609
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03970866666666666
$ws.Range("N4").Value = 0.119126
$ws.Range("Q4").Value = 0.531637776194
$ws.Range("R4").Value = 4.784739985746
$ws.Range("S4").Value = 0.3945718941859609
$ws.Range("T4").Value = 0.3945718941859609
